$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.672.65"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.474.60"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'318.38"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "'92.82"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").Value = "'33.14"
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("D11").Value = "'0.0853"
$ws.Range("E11").Value = "  +8.06%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "2.855.67"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'6.91"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "'15.82"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "2.477.22"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("D18").Value = "41.634.34"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'71.27"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "'11.32"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").Value = "'239.64"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'24.82"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +2.89%  "
$ws.Range("D29").Value = "'9.86"
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("D30").Value = "'36.17"
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("D31").Value = "'159.35"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("D32").Value = "'5.52"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").Value = "'17.32"
$ws.Range("D37").Value = "'1.88"
$ws.Range("E37").Value = "  +5.05%  "
$ws.Range("D38").Value = "'2.93"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").Value = "'0.116"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("D43").Value = "1.993.33"
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.11"
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0285"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("D48").Value = "2.712.34"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "'97.32"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").Value = "'74.15"
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("D51").Value = "'67.08"
$ws.Range("E51").Value = "  +0.01%  "
